# Weekly update: insert a new daily price record for
# "Agrícola del Norte S.A. de Arica - Cebollín baby" as row 43,
# pushing the existing rows 43:82 down to 44:83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43 (shifts old rows 43-82 to 44-83).
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record's data.
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value = "Arica y Parinacota"
$ws.Range("D43").Value = 44673
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = 100112038
$ws.Range("G43").Value = "Cebollín baby"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 250
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 2500
$ws.Range("M43").Value = 2250
$ws.Range("N43").Value = "`$/paquete 1,5 a 2 kilos"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value = 1125
$ws.Range("Q43").Value = 2
$ws.Range("R43").Value = "Hortaliza"
